# Exchange passives (capacitors with higher voltage rating, resistors with
# higher power rating) and refresh the Mouser BOM columns; also un-hide the
# rows that were previously hidden by the autofilter and clear the filter
# criterion (the autofilter range itself stays in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Un-hide every data row first (writing into a hidden row afterwards can
#     stamp a spurious custom row-height, so do this before touching values) --
for ($i = 2; $i -le 23; $i++) {
    $ws.Rows.Item($i).Hidden = $false
}

# --- Column B (Value) updates -------------------------------------------------
$ws.Range("B2").Value  = "0.1uF / 25V"
$ws.Range("B5").Value  = "10uF / 25V"
$ws.Range("B8").Value  = "22u / 25V"
$ws.Range("B10").Value = "6800uF / 16V"
$ws.Range("B12").Value = "MCP73871-2CC"

# --- Column F (Description) updates / additions -------------------------------
$ws.Range("F9").Value  = "Zener Diode"
$ws.Range("F10").Value = "CAPACITOR, European symbol"
$ws.Range("F14").Value = "IC, Ideal Diode Controller"
$ws.Range("F19").Value = "NTC Thermistor"
$ws.Range("F21").Value = "IC, Switching Voltage Regulator"
$ws.Range("F22").Value = "USB-A USB 2.0 Receptacle"
$ws.Range("F23").Value = "Fixed Inductor 1.5uH"

# --- Column G (mouser part no.) updates ---------------------------------------
$ws.Range("G2").Value  = "81-GCM21BR71E104KA7L"
$ws.Range("G3").Value  = "667-ERJ-P06F1003V"
$ws.Range("G4").Value  = "667-ERJ-P06F1002V"
$ws.Range("G5").Value  = "81-GRM21BR61E106MA3L"
$ws.Range("G6").Value  = "667-ERJ-P06F1503V"
$ws.Range("G7").Value  = "667-ERJ-P06F1001V"
$ws.Range("G8").Value  = "81-KRM21FR61E226MF1L"
$ws.Range("G9").Value  = "833-3SMBJ5921B-TP"
$ws.Range("G10").Value = "647-UFW1C682MHD"
$ws.Range("G11").Value = "667-ERJ-P06F4992V"
$ws.Range("G13").Value = "667-ERJ-PB6D7502V"
$ws.Range("G14").Value = "595-LM74670QDGKRQ1"
$ws.Range("G18").Value = "200-MUSBR05FOBSMA"
$ws.Range("G19").Value = "81-NXFT15XH103FEAB45"
$ws.Range("G20").Value = "78-SQ2310ES-T1_BE3"
$ws.Range("G22").Value = "200-USBASSBSM2"

# --- Clear the autofilter criterion but keep the autofilter range -------------
$ws.Range("A1:G23").AutoFilter(6)

# --- Move the active-cell selection, matching the saved view state ------------
$ws.Range("E33").Select()
